$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 129, shifting existing rows 129-134 down to 130-135
$ws.Rows.Item(129).Insert()

# Populate the newly inserted row 129 with the new record's data.
# Columns A, B, C, E, F, G, H, I, R are identical to the neighboring rows for this subset.
$ws.Cells.Item(129, 1).Value = 11
$ws.Cells.Item(129, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(129, 3).Value = "Bíobío"
$ws.Cells.Item(129, 4).Value = 44747
$ws.Cells.Item(129, 4).NumberFormat = $ws.Cells.Item(130, 4).NumberFormat
$ws.Cells.Item(129, 5).Value = 8
$ws.Cells.Item(129, 6).Value = 100112032
$ws.Cells.Item(129, 7).Value = "Zapallo italiano"
$ws.Cells.Item(129, 8).Value = "Sin especificar"
$ws.Cells.Item(129, 9).Value = "Primera"
$ws.Cells.Item(129, 10).Value = 250
$ws.Cells.Item(129, 11).Value = 12000
$ws.Cells.Item(129, 12).Value = 13000
$ws.Cells.Item(129, 13).Value = 12400
$ws.Cells.Item(129, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(129, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(129, 16).Value = 248
$ws.Cells.Item(129, 17).Value = 50
$ws.Cells.Item(129, 18).Value = "Hortaliza"
